$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "42.489.68"
$ws.Cells.Item(2, 5).Value = "  -2.43%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.336.43"
$ws.Cells.Item(3, 5).Value = "  -3.18%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.09%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "319.28"
$ws.Cells.Item(5, 5).Value = "  -2.46%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "103.51"
$ws.Cells.Item(6, 5).Value = "  -0.86%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.636"
$ws.Cells.Item(7, 5).Value = "  -1.38%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.03%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.610"
$ws.Cells.Item(9, 5).Value = "  -6.79%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "39.72"
$ws.Cells.Item(10, 5).Value = "  -6.10%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0917"
$ws.Cells.Item(11, 5).Value = "  -2.74%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.29"
$ws.Cells.Item(12, 5).Value = "  -5.14%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.106"
$ws.Cells.Item(13, 5).Value = "  -0.61%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.976"
$ws.Cells.Item(14, 5).Value = "  -5.61%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.83"
$ws.Cells.Item(15, 5).Value = "  -8.48%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.690.65"
$ws.Cells.Item(16, 5).Value = "  -3.22%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.332.95"
$ws.Cells.Item(17, 5).Value = "  -7.88%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "42.475.13"
$ws.Cells.Item(18, 5).Value = "  -2.51%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.66"
$ws.Cells.Item(19, 5).Value = "  +3.30%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "0.0000105"
$ws.Cells.Item(20, 5).Value = "  -4.36%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "76.13"
$ws.Cells.Item(21, 5).Value = "  +0.62%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "3.56"
$ws.Cells.Item(22, 5).Value = "  +1.06%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "258.48"
$ws.Cells.Item(23, 5).Value = "  -1.05%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.28"
$ws.Cells.Item(24, 5).Value = "  -6.77%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.67"
$ws.Cells.Item(25, 5).Value = "  -0.28%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  +0.07%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "11.26"
$ws.Cells.Item(27, 5).Value = "  -6.74%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "22.93"
$ws.Cells.Item(28, 5).Value = "  -0.26%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -1.20%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "174.05"
$ws.Cells.Item(30, 5).Value = "  -2.57%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "35.17"
$ws.Cells.Item(31, 5).Value = "  -8.66%  "

# Row 32
$ws.Cells.Item(32, 2).Value = "WEMIXToken"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "2.99"
$ws.Cells.Item(32, 5).Value = "  -7.88%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Hedera"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0883"
$ws.Cells.Item(33, 5).Value = "  -5.63%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.99"
$ws.Cells.Item(34, 5).Value = "  -0.19%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  -1.67%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.111"
$ws.Cells.Item(36, 5).Value = "  +2.57%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.50"
$ws.Cells.Item(37, 5).Value = "  -8.81%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0352"
$ws.Cells.Item(38, 5).Value = "  -5.27%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.69"
$ws.Cells.Item(39, 5).Value = "  -6.74%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.60"
$ws.Cells.Item(40, 5).Value = "  -12.08%  "

# Row 41
$ws.Cells.Item(41, 5).Value = "  -11.07%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.229"
$ws.Cells.Item(42, 5).Value = "  -2.48%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "69.28"
$ws.Cells.Item(43, 5).Value = "  -0.29%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  -0.08%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "114.57"
$ws.Cells.Item(45, 5).Value = "  -8.05%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -4.27%  "

# Row 47
$ws.Cells.Item(47, 2).Value = "Celestia"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "11.60"
$ws.Cells.Item(47, 5).Value = "  -8.99%  "

# Row 48
$ws.Cells.Item(48, 2).Value = "FraxShare"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.09"
$ws.Cells.Item(48, 5).Value = "  -4.42%  "

# Row 49
$ws.Cells.Item(49, 2).Value = "BitcoinSV"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "84.57"
$ws.Cells.Item(49, 5).Value = "  +8.78%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "ordi"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "71.92"
$ws.Cells.Item(50, 5).Value = "  +1.15%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "0.0990"
$ws.Cells.Item(51, 5).Value = "  -1.95%  "
